# NIT-9009772691.xlsx — "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker detail table (rows 15-17) listed two workers in arrears:
#   row16: CC  23002892  LUZCELIS ZUÑIGA AYALA      2507  56940  781242
#   row17: CE  334351    ROMAN GABRIEL MONTELLANO    1607  27600  690000
#
# This update removes the first worker's record (LUZCELIS ZUÑIGA AYALA) from
# the statement, leaving only the ROMAN GABRIEL MONTELLANO record (which
# shifts up into row 16), and refreshes the summary figures above the table
# accordingly (total overdue value and the worker/period counters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the worker row for "LUZCELIS ZUÑIGA AYALA" (CC 23002892) entirely;
# Excel shifts every row below it (including the ROMAN GABRIEL MONTELLANO
# row and the signature block) up by one.
$ws.Rows(16).EntireRow.Delete()

# Refresh the summary header figures to reflect the single remaining record.
$ws.Range("E11").Value = 27600   # VALOR MORA total
$ws.Range("C13").Value = 1       # Cant. Trabajadores
$ws.Range("F13").Value = 1       # Cant. Periodos
